$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'42.824.70"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = "'2.299.47"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'305.96"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.84%  '
$ws.Range('D6').Value = "'96.85"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.75%  '
$ws.Range('E7').Value = '  -1.72%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -2.58%  '
$ws.Range('D10').Value = "'35.55"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.00%  '
$ws.Range('E11').Value = '  +0.28%  '
$ws.Range('D12').Value = "'18.37"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.62%  '
$ws.Range('E13').Value = '  +1.15%  '
$ws.Range('E14').Value = '  -1.88%  '
$ws.Range('D15').Value = "'2.658.35"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.03%  '
$ws.Range('D16').Value = "'2.292.73"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.07%  '
$ws.Range('D17').Value = "'0.781"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.87%  '
$ws.Range('D18').Value = "'42.769.62"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.46%  '
$ws.Range('D19').Value = "'13.05"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('D20').Value = "'0.0₃0899"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.18%  '
$ws.Range('E21').Value = '  -1.73%  '
$ws.Range('D22').Value = "'67.37"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.25%  '
$ws.Range('D23').Value = "'236.17"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.81%  '
$ws.Range('D24').Value = "'2.14"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.68%  '
$ws.Range('E25').Value = '  +1.70%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').Value = "'25.40"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.49%  '
$ws.Range('D29').Value = "'166.24"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.23%  '
$ws.Range('E30').Value = '  +1.41%  '
$ws.Range('E31').Value = '  -0.90%  '
$ws.Range('D32').Value = "'33.33"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.69%  '
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('D34').Value = "'4.76"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.46%  '
$ws.Range('E35').Value = '  -2.55%  '
$ws.Range('D36').Value = "'17.74"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.18%  '
$ws.Range('D38').Value = "'0.0692"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('E39').Value = '  -0.87%  '
$ws.Range('E40').Value = '  -1.31%  '
$ws.Range('E41').Value = '  -1.23%  '
$ws.Range('E42').Value = '  -1.67%  '
$ws.Range('D43').Value = "'2.003.06"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.36%  '
$ws.Range('D44').Value = "'0.0281"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.32%  '
$ws.Range('D45').Value = "'18.25"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.26%  '
$ws.Range('D46').Value = "'9.98"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.30%  '
$ws.Range('E47').Value = '  -5.16%  '
$ws.Range('E48').Value = '  -2.63%  '
$ws.Range('E49').Value = '  +7.11%  '
$ws.Range('D50').Value = "'53.69"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.12%  '
$ws.Range('D51').Value = "'2.526.73"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.39%  '
